$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 to use the new "Опис багу" string value
$ws.Range("B1").Value = "Опис багу"

# Move the selection/view to B6 (also clears the stale topLeftCell scroll anchor)
$ws.Range("B6").Select()
